# #5: fund, bonds, otherbonds, antique done
#
# Rework the "其他有價證券" sheet (sheet 6) into "具有相當價值之財產":
#  - rename the worksheet
#  - collapse the old 4-row / 7-column layout into the standard
#    1 header row + 1 data row / 12 column layout used by the other sheets
#  - row 2 now holds the golf-club-membership-card record

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- rename sheet ---------------------------------------------------------
$ws.Name = "具有相當價值之財產"

# --- drop the two extra legacy rows (old header/template rows) -----------
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# --- extend header row (row 1) and data row (row 2) out to column L,
#     copying formatting from existing cells so the new cells pick up the
#     same style index (s="1" for row 1, s="2"/"1" for row 2) -------------
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("G1").Copy($ws.Range("I1"))
$ws.Range("G1").Copy($ws.Range("J1"))
$ws.Range("G1").Copy($ws.Range("K1"))
$ws.Range("G1").Copy($ws.Range("L1"))

$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Range("G2").Copy($ws.Range("I2"))
$ws.Range("G2").Copy($ws.Range("J2"))
$ws.Range("G2").Copy($ws.Range("K2"))
$ws.Range("A2").Copy($ws.Range("L2"))

# --- header row values ------------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "property_category"
$ws.Range("G1").Value = "category"
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("K1").Value = "source_file"
$ws.Range("L1").Value = "index"

# --- data row values (golf club membership card, owned by 蔡慧敏) --------
$ws.Range("A2").Value = 87
$ws.Range("B2").Value = "統帥球場高爾夫球証"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "蔡慧敏"
$ws.Range("E2").Value = 500000
$ws.Range("F2").Value = "otherbonds"
$ws.Range("G2").Value = "normal"
$ws.Range("H2").Value = "2011-11-22"
$ws.Range("I2").Value = "陳根德"
$ws.Range("J2").Value = 833
$ws.Range("K2").Value = "tmpa3b61"
$ws.Range("L2").Value = 87
